# Insert a new data row at row 515 (pushes the existing rows 515-629 down to 516-630)
# and populate it with the new record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("515:515").Insert()

$ws.Range("A515").Value = 10
$ws.Range("B515").Value = "Vega Modelo de Temuco"
$ws.Range("C515").Value = "La Araucanía"
$ws.Range("D515").Value = 44889
$ws.Range("E515").Value = 9
$ws.Range("F515").Value = 100112043
$ws.Range("G515").Value = "Pepino ensalada"
$ws.Range("H515").Value = "Sin especificar"
$ws.Range("I515").Value = "Primera"
$ws.Range("J515").Value = 125
$ws.Range("K515").Value = 21000
$ws.Range("L515").Value = 21000
$ws.Range("M515").Value = 21000
$ws.Range("N515").Value = "$/caja 60 unidades"
$ws.Range("O515").Value = "Región de O'Higgins"
$ws.Range("P515").Value = 350
$ws.Range("Q515").Value = 60
$ws.Range("R515").Value = "Hortaliza"
